$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question 3 : ajout du facteur PSF (C1 = explication, C2 = calcul)
$ws.Range("C1").Value = "Pour intégrer cette information dans notre calcul, nous ajoutons un facteur de multiplication (PSF). Dans notre cas, ce facteur correspond au nombre de niveaux d'énergie, nous multiplions donc le résultat obtenu par 11."
$ws.Range("C2").Formula = "=(B1*11)"

# Reprend la mise en forme (fond jaune) de B1 pour la nouvelle cellule de résultat C2
$ws.Range("B1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Met à jour la sélection active comme dans le classeur d'origine
$ws.Range("C2").Select()
